# Fix typo in presentation: slide 14, body placeholder, 2nd bullet.
# "Намаляне на размерността на пространството - ..." ->
# "Намаляване на размерността на пространството - ..."
#
# The original paragraph is split into 3 runs:
#   r1 = "Н"
#   r2 = "амаляне на размерността на пространството - "
#   r3 = "Singular Value Decomposition "
#
# The target paragraph is split into 4 runs:
#   r1 = "Намаляване"
#   r2 = " "
#   r3 = "на размерността на пространството - "
#   r4 = "Singular Value Decomposition "   (unchanged)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)
$shape = $s.Shapes.Item(2)
$textRange = $shape.TextFrame.TextRange
$para = $textRange.Paragraphs(2)

# 1) Fix the typo: expand the first run "Н" into "Намаляване".
$para.Characters(1, 1).Text = "Намаляване"

# 2) Drop the now-redundant "амаляне" prefix that used to follow "Н" in the
#    second run, leaving the pre-existing space before "на ...".
$para.Characters(11, 7).Text = ""

# 3) Re-touch the leading space on its own so it keeps being its own run,
#    separate from the "на размерността ..." text that follows it.
$para.Characters(11, 1).Text = " "
